# Update the "想去人数" (interest count) column (F) with freshly scraped
# numbers on both the "展览" sheet and the aggregated "全部类型" sheet.
# (gh-pages data refresh -> output generated at 456a3b4)

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 34
$ws1.Range("F3").Value = 772
$ws1.Range("F6").Value = 50
$ws1.Range("F7").Value = 267
$ws1.Range("F8").Value = 3753
$ws1.Range("F9").Value = 75
$ws1.Range("F10").Value = 4417
$ws1.Range("F12").Value = 1098
$ws1.Range("F13").Value = 59

# Sheet "全部类型" (all types) - same underlying events, shifted row offset
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 34
$ws4.Range("F3").Value = 772
$ws4.Range("F6").Value = 50
$ws4.Range("F8").Value = 267
$ws4.Range("F9").Value = 3753
$ws4.Range("F10").Value = 75
$ws4.Range("F11").Value = 4417
$ws4.Range("F13").Value = 1098
$ws4.Range("F14").Value = 59
